$d = $word.ActiveDocument

# 1) Update the letter date: September 19, 2025 -> September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2) Split the mailing address line into two lines and add a trailing blank
#    line: "2443 Fillmore St, San Francisco CA 94301-1055" ->
#      "2443 Fillmore St"
#      "San Francisco, CA 94301-1055"
#      <blank paragraph>
$rng = $d.Content
$rng.Find.Execute(", San Francisco CA 94301-1055", $true, $false, $false, $false, $false, $true, 1, $false, "^pSan Francisco, CA 94301-1055^p", 2) | Out-Null

# 3) Remove the two empty paragraphs that used to sit right after the
#    "Board of Directors" signature line.
$idx = 0
$target = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Board of Directors*") {
        $target = $idx
    }
}

$pFirst = $d.Paragraphs($target + 1)
$pSecond = $d.Paragraphs($target + 2)
$pSecond.Range.Delete() | Out-Null
$pFirst.Range.Delete() | Out-Null
